$d = $word.ActiveDocument
$insertPoint = $d.Content
$insertPoint.Collapse(0)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr><w:i w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr><w:i w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr><w:i w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr><w:i w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr><w:i w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr><w:i w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr><w:i w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:rPr><w:b w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Fecha 19-12-2018</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">1.-</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> Cambios en la clase Horario, se agrega el constructor a la clase, el cual crea una semana con 6 días y cada día con 10 bloques. Ahora se pueden asignar cursos a bloques de cualquier día dentro de la planificación semanal (6 días con 10 bloques cada uno (desde las 8:30 en adelante)) de un curso, sala o profesor.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">se agregaron los siguientes métodos: </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">public void visualizarPlanificacionSemanal(){}</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">public boolean verificarDisponibilidadDeBloque(int dia, int bloque){}</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">public void asignarCursoABloque(Curso curso, int dia, int bloque){</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">public void quitarCursoDeBloque(int dia, int bloque){}</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="720"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">2.-</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> Se agregó un método a la clase Curso.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:tab/><w:t xml:space="preserve"> public void visualizarDatos(){}</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">3.-</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Cambios en la clase Gestion, se agrego un método para visualizar los cursos contenidos en la lista.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> </w:t><w:tab/><w:t xml:space="preserve">public void mostrarListaDeCursos(){</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">4</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve"> Se crearon instancias de Curso, permitiendo crear un curso con un nombre y asignarle/quitarle un profesor., se añadieron los sig métodos:</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:tab/><w:t xml:space="preserve">public void deleteProfesor() {}</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">5.- </w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Cambios en el constructor de la clase Sala.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr/></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">6.-</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Se crearon instancias de Sala, permitiendo crear una sala con un nombre, para posteriormente asignarle un curso a un determinado bloque. También se puede consultar por la disponibilidad de un determinado  bloque de un determinado día. </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720" w:firstLine="0"/><w:rPr><w:i w:val="1"/></w:rPr></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($xml)
